{"js": "// Update each two-digit multiplication equation in the table to the\n// newly generated values (commit \"Update master to output generated at\n// c986bee\"). Each old equation string is unique in the document, so we\n// look it up with body.search() (exact, case-sensitive, whole match)\n// and replace the text of the single hit.\n\nconst replacements = [\n  [\"35\u00d772=2520\", \"90\u00d742=3780\"],\n  [\"84\u00d748=4032\", \"95\u00d773=6935\"],\n  [\"71\u00d737=2627\", \"93\u00d752=4836\"],\n  [\"11\u00d742=462\", \"55\u00d786=4730\"],\n  [\"68\u00d750=3400\", \"94\u00d745=4230\"],\n  [\"33\u00d718=594\", \"59\u00d785=5015\"],\n  [\"12\u00d730=360\", \"64\u00d760=3840\"],\n  [\"89\u00d750=4450\", \"42\u00d735=1470\"],\n  [\"15\u00d768=1020\", \"38\u00d739=1482\"],\n  [\"99\u00d762=6138\", \"91\u00d715=1365\"],\n  [\"66\u00d739=2574\", \"59\u00d786=5074\"],\n  [\"60\u00d724=1440\", \"82\u00d738=3116\"],\n  [\"28\u00d776=2128\", \"70\u00d749=3430\"],\n  [\"44\u00d784=3696\", \"33\u00d771=2343\"],\n  [\"90\u00d723=2070\", \"14\u00d798=1372\"],\n  [\"19\u00d735=665\", \"70\u00d740=2800\"],\n  [\"29\u00d748=1392\", \"64\u00d727=1728\"],\n  [\"17\u00d738=646\", \"35\u00d772=2520\"],\n  [\"45\u00d716=720\", \"19\u00d765=1235\"],\n  [\"28\u00d726=728\", \"12\u00d756=672\"],\n  [\"74\u00d720=1480\", \"79\u00d729=2291\"],\n  [\"86\u00d780=6880\", \"21\u00d778=1638\"],\n  [\"58\u00d787=5046\", \"43\u00d785=3655\"],\n  [\"36\u00d711=396\", \"42\u00d731=1302\"],\n  [\"14\u00d748=672\", \"61\u00d720=1220\"],\n];\n\nconst body = context.document.body;\n\n// Process replacements one at a time (sync between each) so that a new\n// value that happens to equal an old value elsewhere in the list (e.g.\n// \"35\u00d772=2520\" is both an original value and a later replacement's new\n// value) never gets matched again by a subsequent search \u2014 each search\n// runs only after the previous replacement has already been committed.\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  // Only the first occurrence should exist (the strings are unique),\n  // but guard with a loop in case of duplicates.\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update each two-digit multiplication equation in the table to the\n# newly generated values (commit \"Update master to output generated at\n# c986bee\"). Every old equation string is unique in the document, so a\n# simple Find/Replace (one hit, ReplaceOne) on the whole document body\n# locates and rewrites each cell's text in turn.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"35\u00d772=2520\", \"90\u00d742=3780\"),\n    @(\"84\u00d748=4032\", \"95\u00d773=6935\"),\n    @(\"71\u00d737=2627\", \"93\u00d752=4836\"),\n    @(\"11\u00d742=462\",  \"55\u00d786=4730\"),\n    @(\"68\u00d750=3400\", \"94\u00d745=4230\"),\n    @(\"33\u00d718=594\",  \"59\u00d785=5015\"),\n    @(\"12\u00d730=360\",  \"64\u00d760=3840\"),\n    @(\"89\u00d750=4450\", \"42\u00d735=1470\"),\n    @(\"15\u00d768=1020\", \"38\u00d739=1482\"),\n    @(\"99\u00d762=6138\", \"91\u00d715=1365\"),\n    @(\"66\u00d739=2574\", \"59\u00d786=5074\"),\n    @(\"60\u00d724=1440\", \"82\u00d738=3116\"),\n    @(\"28\u00d776=2128\", \"70\u00d749=3430\"),\n    @(\"44\u00d784=3696\", \"33\u00d771=2343\"),\n    @(\"90\u00d723=2070\", \"14\u00d798=1372\"),\n    @(\"19\u00d735=665\",  \"70\u00d740=2800\"),\n    @(\"29\u00d748=1392\", \"64\u00d727=1728\"),\n    @(\"17\u00d738=646\",  \"35\u00d772=2520\"),\n    @(\"45\u00d716=720\",  \"19\u00d765=1235\"),\n    @(\"28\u00d726=728\",  \"12\u00d756=672\"),\n    @(\"74\u00d720=1480\", \"79\u00d729=2291\"),\n    @(\"86\u00d780=6880\", \"21\u00d778=1638\"),\n    @(\"58\u00d787=5046\", \"43\u00d785=3655\"),\n    @(\"36\u00d711=396\",  \"42\u00d731=1302\"),\n    @(\"14\u00d748=672\",  \"61\u00d720=1220\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #         MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n    #         Format, ReplaceWith, Replace)\n    # Wrap=1 -> wdFindContinue, Replace=2 -> wdReplaceOne (replace the\n    # single match this unique string has, then stop).\n    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"WARNING: could not find '$oldText'\"\n    }\n}\n"}
